$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coefficient/constant values for rows 11, 12 and 17
# Row 11: SteelOxygenBlownConverter
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 2

# Row 12: SteelElectricFurnace
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 2

# Row 17: PigIron
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = 2

# Update the active selection to B17
$ws.Range("B17").Select()
